$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in the header (F1)
$ws.Range("F1").Value = "Last status check on: 08.02.2022 09:30"

# Row 6 (Shell Olomoucká) gets a fresh price reading:
#  - new current price in B6, previous price shifts into C6
#  - delta (D6) and old-datum (E6) are now reported as plain text values
#    straight from the scraping script instead of numbers/dates
$ws.Range("B6").Value = 37.5
$ws.Range("C6").Value = 36.9

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "+0.6"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2022-02-08 09:30:21"
$ws.Range("E6").Style = "Normal"
